{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of the change (per the OOXML diff):\n//   1. Heading \"4.2.12\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\" -> \"4.2.12\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\"\n//      (and its bookmark is renamed to match: \u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599 -> \u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599)\n//   2. Heading \"4.2.19\u9000\u51fa\u7cfb\u7edf\" -> \"4.2.19\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\"\n//      (its bookmark \"\u9000\u51fa\u7cfb\u7edf-1\" is renamed to \"\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\")\n//      Its body paragraph text becomes the new \"\u5f53\u7ba1\u7406\u5458\u8f93\u5165\u7528\u6237\u8d44\u6599...\" text.\n//   3. A new \"4.2.20\u9000\u51fa\u7cfb\u7edf\" Heading4 paragraph (with a fresh bookmark named\n//      \"\u9000\u51fa\u7cfb\u7edf-1\") is inserted right after that, followed by a paragraph\n//      that holds the original \"\u5f53\u7cfb\u7edf\u7ba1\u7406\u5458\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e...\" text (i.e. the\n//      old body paragraph's content now lives one slot further down).\n\nconst doc = context.document;\n\n// ---------------------------------------------------------------------\n// Part 1: \"\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\" -> \"\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\"\n// ---------------------------------------------------------------------\nconst bm1 = doc.getBookmarkRange(\"\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\");\nbm1.load(\"text\");\nawait context.sync();\n\n// The bookmark is a zero-length marker that sits right before the Heading4\n// run, so its paragraph is the heading paragraph we need to retext.\nconst heading1 = bm1.paragraphs.getFirst();\nheading1.load(\"text\");\nawait context.sync();\n\n// Rename the bookmark (delete + reinsert at the same collapsed position).\ndoc.deleteBookmark(\"\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\");\nbm1.insertBookmark(\"\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\");\n\n// Update the heading text itself.\nheading1.getRange().insertText(\"4.2.12\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\", Word.InsertLocation.replace);\n\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Part 2: \"4.2.19\u9000\u51fa\u7cfb\u7edf\" -> \"4.2.19\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\" (+ new paragraph content)\n//         and insertion of a new \"4.2.20\u9000\u51fa\u7cfb\u7edf\" section after it.\n// ---------------------------------------------------------------------\nconst bm2 = doc.getBookmarkRange(\"\u9000\u51fa\u7cfb\u7edf-1\");\nbm2.load(\"text\");\nawait context.sync();\n\nconst heading2 = bm2.paragraphs.getFirst();  // \"4.2.19\u9000\u51fa\u7cfb\u7edf\"\nheading2.load(\"text\");\nawait context.sync();\nconst body2 = heading2.getNext();            // \"\u5f53\u7cfb\u7edf\u7ba1\u7406\u5458\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e...\"\nbody2.load(\"text\");\nawait context.sync();\n\n// Rename the bookmark \"\u9000\u51fa\u7cfb\u7edf-1\" -> \"\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\" on the (still) 4.2.19 heading.\ndoc.deleteBookmark(\"\u9000\u51fa\u7cfb\u7edf-1\");\nbm2.insertBookmark(\"\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\");\n\n// Update heading + body paragraph texts.\nheading2.getRange().insertText(\"4.2.19\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\", Word.InsertLocation.replace);\nbody2.getRange().insertText(\n  \"\u5f53\u7ba1\u7406\u5458\u8f93\u5165\u7528\u6237\u8d44\u6599\u3001\u4fee\u6539\u540e\u9700\u8981\u7ba1\u7406\u5458\u51b3\u5b9a\u662f\u5426\u4fdd\u5b58\u7528\u6237\u8d44\u6599\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n\n// Insert the new \"4.2.20\u9000\u51fa\u7cfb\u7edf\" heading paragraph right after the\n// (retexted) body paragraph, followed by a paragraph holding the original\n// \"\u5f53\u7cfb\u7edf\u7ba1\u7406\u5458\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e...\" text.\nconst newHeading = body2.insertParagraph(\"4.2.20\u9000\u51fa\u7cfb\u7edf\", Word.InsertLocation.after);\nnewHeading.style = \"Heading 4\";\nawait context.sync();\n\nconst newBody = newHeading.insertParagraph(\n  \"\u5f53\u7cfb\u7edf\u7ba1\u7406\u5458\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e\uff0c\u53ef\u9000\u51fa\u8be5\u7cfb\u7edf\u3002\u5982\u679c\u7528\u6237\u8fdb\u884c\u4e86\u5f71\u54cd\u7528\u6237\u4fe1\u606f\u7684\u64cd\u4f5c\uff0c\u63d0\u793a\u7cfb\u7edf\u7ba1\u7406\u5458\u662f\u5426\u8fdb\u884c\u4fdd\u5b58\u3002\",\n  Word.InsertLocation.after\n);\nnewBody.style = \"First Paragraph\";\nawait context.sync();\n\n// Re-fetch the new heading paragraph through the body's paragraph collection\n// (a freshly loaded collection, not the proxy returned by insertParagraph)\n// before attaching the bookmark, so the bookmarkStart/bookmarkEnd pair lands\n// correctly within the same paragraph.\nconst allParas = doc.body.paragraphs;\nallParas.load(\"text\");\nawait context.sync();\n\nlet newHeadingIndex = -1;\nfor (let i = 0; i < allParas.items.length; i++) {\n  if (allParas.items[i].text === \"4.2.20\u9000\u51fa\u7cfb\u7edf\") {\n    newHeadingIndex = i;\n    break;\n  }\n}\nconst freshNewHeading = allParas.items[newHeadingIndex];\nfreshNewHeading.getRange(\"Start\").insertBookmark(\"\u9000\u51fa\u7cfb\u7edf-1\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Summary of the change (per the OOXML diff):\n#   1. Heading \"4.2.12\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\" -> \"4.2.12\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\"\n#      (and its bookmark is renamed to match: \u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599 -> \u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599)\n#   2. Heading \"4.2.19\u9000\u51fa\u7cfb\u7edf\" -> \"4.2.19\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\"\n#      (its bookmark \"\u9000\u51fa\u7cfb\u7edf-1\" is renamed to \"\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\")\n#      Its body paragraph text becomes the new \"\u5f53\u7ba1\u7406\u5458\u8f93\u5165\u7528\u6237\u8d44\u6599...\" text.\n#   3. A new \"4.2.20\u9000\u51fa\u7cfb\u7edf\" Heading4 paragraph (with a fresh bookmark named\n#      \"\u9000\u51fa\u7cfb\u7edf-1\") is inserted right after that, followed by a paragraph\n#      that holds the original \"\u5f53\u7cfb\u7edf\u7ba1\u7406\u5458\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e...\" text (i.e. the\n#      old body paragraph's content now lives one slot further down).\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Part 1: \"\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\" -> \"\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\"\n# ---------------------------------------------------------------------\n$bm1 = $d.Bookmarks(\"\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\")\n$bm1Range = $bm1.Range\n\n# The bookmark is a zero-length marker that sits right before the Heading4\n# run, so its paragraph is the heading paragraph we need to retext.\n$heading1 = $bm1Range.Paragraphs(1)\n\n# Rename the bookmark (delete + re-add at the same collapsed position).\n$bm1.Delete()\n$d.Bookmarks.Add(\"\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\", $bm1Range)\n\n# Update the heading text itself.\n$heading1.Range.Text = \"4.2.12\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u5de5\u7a0b\u5e08\u8d44\u6599\"\n\n# ---------------------------------------------------------------------\n# Part 2: \"4.2.19\u9000\u51fa\u7cfb\u7edf\" -> \"4.2.19\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\" (+ new paragraph content)\n#         and insertion of a new \"4.2.20\u9000\u51fa\u7cfb\u7edf\" section after it.\n# ---------------------------------------------------------------------\n$bm2 = $d.Bookmarks(\"\u9000\u51fa\u7cfb\u7edf-1\")\n$bm2Range = $bm2.Range\n$heading2 = $bm2Range.Paragraphs(1)      # \"4.2.19\u9000\u51fa\u7cfb\u7edf\"\n$body2 = $heading2.Next()                 # \"\u5f53\u7cfb\u7edf\u7ba1\u7406\u5458\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e...\"\n\n# Rename the bookmark \"\u9000\u51fa\u7cfb\u7edf-1\" -> \"\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\" on the (still) 4.2.19 heading.\n$bm2.Delete()\n$d.Bookmarks.Add(\"\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\", $bm2Range)\n\n# Update heading + body paragraph texts.\n$heading2.Range.Text = \"4.2.19\u4fdd\u5b58\u7528\u6237\u4fe1\u606f\"\n$body2.Range.Text = \"\u5f53\u7ba1\u7406\u5458\u8f93\u5165\u7528\u6237\u8d44\u6599\u3001\u4fee\u6539\u540e\u9700\u8981\u7ba1\u7406\u5458\u51b3\u5b9a\u662f\u5426\u4fdd\u5b58\u7528\u6237\u8d44\u6599\"\n\n# Insert the new \"4.2.20\u9000\u51fa\u7cfb\u7edf\" heading paragraph right after the\n# (retexted) body paragraph, followed by a paragraph holding the original\n# \"\u5f53\u7cfb\u7edf\u7ba1\u7406\u5458\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e...\" text.\n$body2.Range.InsertParagraphAfter()\n$newHeading = $body2.Next()\n$newHeading.Range.Text = \"4.2.20\u9000\u51fa\u7cfb\u7edf\"\n$newHeading.Style = \"Heading 4\"\n\n$newHeading.Range.InsertParagraphAfter()\n$newBody = $newHeading.Next()\n$newBody.Range.Text = \"\u5f53\u7cfb\u7edf\u7ba1\u7406\u5458\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e\uff0c\u53ef\u9000\u51fa\u8be5\u7cfb\u7edf\u3002\u5982\u679c\u7528\u6237\u8fdb\u884c\u4e86\u5f71\u54cd\u7528\u6237\u4fe1\u606f\u7684\u64cd\u4f5c\uff0c\u63d0\u793a\u7cfb\u7edf\u7ba1\u7406\u5458\u662f\u5426\u8fdb\u884c\u4fdd\u5b58\u3002\"\n$newBody.Style = \"First Paragraph\"\n\n# Attach the \"\u9000\u51fa\u7cfb\u7edf-1\" bookmark to the new heading paragraph. Collapse the\n# range to its start first so bookmarkStart/bookmarkEnd land adjacent to each\n# other, right before the run, instead of the pair straddling into the next\n# paragraph.\n$newHeadingStart = $newHeading.Range\n$newHeadingStart.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"\u9000\u51fa\u7cfb\u7edf-1\", $newHeadingStart)\n"}
